$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MIT-FS.Audit4Improve-API")

# --- 1. Make room: insert 6 new rows starting at row 10 (old rows 10.. shift down) ---
$ws.Rows.Item(10).Resize(6).Insert()

# --- 2. Refresh the "saved" timestamp header (B1) ---
$ws.Range("B1").Value = "Fri Mar 08 15:38:21 CET 2024"

# --- 3. Rewrite the metric rows 2-15 (name, value, unit, description, source, timestamp) ---
$rows = @(
  @{ r=2;  A="forks";          B="43";                              C="forks";          D="Número de forks, no son los forks de la web";                               E="GitHub";            F="Fri Mar 08 15:38:18 CET 2024" },
  @{ r=3;  A="closedIssues";   B="47";                              C="issues";         D="Numero de asuntos cerrados";                                                E="GitHub";            F="Fri Mar 08 15:38:20 CET 2024" },
  @{ r=4;  A="ownerCommits";   B="52";                              C="commits";        D="Commits del responsable";                                                   E="GitHub";            F="Fri Mar 08 15:38:18 CET 2024" },
  @{ r=5;  A="totalAdditions"; B="6909";                            C="additions";      D="Suma el total de adiciones desde que el repositorio se creó";              E="GitHub, calculada"; F="Fri Mar 08 15:38:18 CET 2024" },
  @{ r=6;  A="totalDeletions"; B="6909";                            C="deletions";      D="Suma el total de eliminaciones desde que el repositorio se cre�";          E="GitHub, calculada"; F="Fri Mar 08 15:38:18 CET 2024" },
  @{ r=7;  A="subscribers";    B="2";                               C="subscribers";    D="Número de suscriptores de un repositorio, watchers en la web";             E="GitHub";            F="Fri Mar 08 15:38:18 CET 2024" },
  @{ r=8;  A="watchers";       B="0";                               C="watchers";       D="Observadores de un repositorio, en la web aparece com forks";              E="GitHub";            F="Fri Mar 08 15:38:18 CET 2024" },
  @{ r=9;  A="stars";          B="0";                               C="stars";          D="Numero de estrellas";                                                       E="GitHub";            F="Fri Mar 08 15:38:18 CET 2024" },
  @{ r=10; A="issues";         B="114";                             C="issues";         D="Numero de asuntos totales";                                                 E="GitHub";            F="Fri Mar 08 15:38:20 CET 2024" },
  @{ r=11; A="openIssues";     B="67";                              C="issues";         D="Numero de asuntos abiertos";                                                E="GitHub";            F="Fri Mar 08 15:38:20 CET 2024" },
  @{ r=12; A="lastPush";       B="Fri Mar 08 13:45:32 CET 2024";    C="date";           D="Último push realizado en el repositorio";                                   E="GitHub";            F="Fri Mar 08 15:38:21 CET 2024" },
  @{ r=13; A="lastUpdated";    B="Tue Apr 18 13:27:41 CEST 2023";   C="date";           D="Última actualización";                                                      E="GitHub";            F="Fri Mar 08 15:38:21 CET 2024" },
  @{ r=14; A="collaborators";  B="28";                              C="collaborators";  D="Numero de colaboradores en el repositorio";                                 E="GitHub";            F="Fri Mar 08 15:38:18 CET 2024" },
  @{ r=15; A="creation";       B="Thu Feb 03 11:04:44 CET 2022";    C="date";           D="Fecha de creación del repositorio";                                         E="GitHub";            F="Fri Mar 08 15:38:21 CET 2024" }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F

    $fmt = $ws.Range("A$r" + ":B$r")
    $fmt.Font.Name = "Serif"
    $fmt.Font.Size = 10
    $fmt.Font.Color = 65280
    $fmt.Font.Bold = $true
    $fmt.Interior.Color = 65280
}

# --- 4. Move the trailing "Indicadores" label down to row 16 ---
$ws.Range("A16").Value = "Indicadores"

# --- 5. Re-fit columns to the new content ---
$ws.Columns.Item(3).EntireColumn.AutoFit()
$ws.Columns.Item(4).EntireColumn.AutoFit()
$ws.Columns.Item(5).EntireColumn.AutoFit()
